$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(161, 2).Value = 53925
$ws.Cells.Item(161, 5).Value = 79.37
$ws.Cells.Item(161, 6).Value = 1
$ws.Cells.Item(161, 7).Value = 66.44
$ws.Cells.Item(162, 2).Value = 64350
$ws.Cells.Item(162, 5).Value = 70.63
$ws.Cells.Item(162, 6).Value = 101
$ws.Cells.Item(162, 7).Value = 6710.44
$ws.Cells.Item(163, 2).Value = 57756
$ws.Cells.Item(163, 6).Value = -100
$ws.Cells.Item(163, 7).Value = -6644
$ws.Cells.Item(264, 2).Value = 64979
$ws.Cells.Item(264, 5).Value = 314.41
$ws.Cells.Item(264, 6).Value = 82
$ws.Cells.Item(264, 7).Value = 24251.5
$ws.Cells.Item(265, 2).Value = 48719
$ws.Cells.Item(265, 5).Value = 353.35
$ws.Cells.Item(265, 6).Value = -81
$ws.Cells.Item(265, 7).Value = -23955.75
$ws.Cells.Item(313, 2).Value = 62997
$ws.Cells.Item(313, 6).Value = 72
$ws.Cells.Item(313, 7).Value = 22020.48
$ws.Cells.Item(314, 2).Value = 57854
$ws.Cells.Item(314, 6).Value = 2
$ws.Cells.Item(314, 7).Value = 611.6799999999999
$ws.Cells.Item(317, 2).Value = 61610
$ws.Cells.Item(317, 4).Value = 102.71
$ws.Cells.Item(317, 5).Value = 122.71
$ws.Cells.Item(317, 6).Value = -58
$ws.Cells.Item(317, 7).Value = -5957.18
$ws.Cells.Item(318, 2).Value = 57077
$ws.Cells.Item(318, 4).Value = 93.08
$ws.Cells.Item(318, 5).Value = 111.2
$ws.Cells.Item(318, 6).Value = 1
$ws.Cells.Item(318, 7).Value = 93.08
$ws.Cells.Item(346, 2).Value = 55373
$ws.Cells.Item(346, 5).Value = 163.62
$ws.Cells.Item(346, 6).Value = -94
$ws.Cells.Item(346, 7).Value = -13562.32
$ws.Cells.Item(347, 2).Value = 63520
$ws.Cells.Item(347, 5).Value = 153.4
$ws.Cells.Item(347, 6).Value = 97
$ws.Cells.Item(347, 7).Value = 13995.16
$ws.Cells.Item(350, 2).Value = 63571
$ws.Cells.Item(350, 5).Value = 152.53
$ws.Cells.Item(350, 6).Value = 27
$ws.Cells.Item(350, 7).Value = 3873.96
$ws.Cells.Item(351, 2).Value = 63531
$ws.Cells.Item(351, 6).Value = 80
$ws.Cells.Item(351, 7).Value = 11478.4
$ws.Cells.Item(352, 2).Value = 57802
$ws.Cells.Item(352, 5).Value = 162.71
$ws.Cells.Item(352, 6).Value = -79
$ws.Cells.Item(352, 7).Value = -11334.92
$ws.Cells.Item(355, 2).Value = 55356
$ws.Cells.Item(355, 5).Value = 54.04
$ws.Cells.Item(355, 6).Value = -158
$ws.Cells.Item(355, 7).Value = -7527.12
$ws.Cells.Item(356, 2).Value = 63510
$ws.Cells.Item(356, 5).Value = 50.66
$ws.Cells.Item(356, 6).Value = 167
$ws.Cells.Item(356, 7).Value = 7955.88
$ws.Cells.Item(372, 2).Value = 57885
$ws.Cells.Item(372, 5).Value = 62.28
$ws.Cells.Item(372, 6).Value = 4
$ws.Cells.Item(372, 7).Value = 208.52
$ws.Cells.Item(373, 2).Value = 63652
$ws.Cells.Item(373, 5).Value = 55.42
$ws.Cells.Item(373, 6).Value = 250
$ws.Cells.Item(373, 7).Value = 13032.5
$ws.Cells.Item(375, 2).Value = 61605
$ws.Cells.Item(375, 5).Value = 133.78
$ws.Cells.Item(375, 6).Value = -13
$ws.Cells.Item(375, 7).Value = -1455.48
$ws.Cells.Item(376, 2).Value = 63563
$ws.Cells.Item(376, 5).Value = 119.04
$ws.Cells.Item(376, 6).Value = 15
$ws.Cells.Item(376, 7).Value = 1679.4
$ws.Cells.Item(389, 2).Value = 62865
$ws.Cells.Item(389, 6).Value = 151
$ws.Cells.Item(389, 7).Value = 12051.31
$ws.Cells.Item(390, 2).Value = 57817
$ws.Cells.Item(390, 6).Value = 3
$ws.Cells.Item(390, 7).Value = 239.43
$ws.Cells.Item(419, 2).Value = 63007
$ws.Cells.Item(419, 6).Value = 984
$ws.Cells.Item(419, 7).Value = 168588.72
$ws.Cells.Item(420, 2).Value = 57856
$ws.Cells.Item(420, 6).Value = 2
$ws.Cells.Item(420, 7).Value = 342.66
$ws.Cells.Item(421, 2).Value = 63008
$ws.Cells.Item(421, 6).Value = 504
$ws.Cells.Item(421, 7).Value = 76189.67999999999
$ws.Cells.Item(422, 2).Value = 57857
$ws.Cells.Item(422, 6).Value = 3
$ws.Cells.Item(422, 7).Value = 453.51
$ws.Cells.Item(431, 2).Value = 53082
$ws.Cells.Item(431, 3).Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Cells.Item(431, 6).Value = 1
$ws.Cells.Item(431, 7).Value = 59.47
$ws.Cells.Item(432, 2).Value = 63102
$ws.Cells.Item(432, 3).Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Cells.Item(432, 6).Value = 36
$ws.Cells.Item(432, 7).Value = 2140.92
$ws.Cells.Item(579, 2).Value = 65069
$ws.Cells.Item(579, 5).Value = 14.3
$ws.Cells.Item(579, 6).Value = 172
$ws.Cells.Item(579, 7).Value = 2313.4
$ws.Cells.Item(580, 2).Value = 53757
$ws.Cells.Item(580, 5).Value = 16.08
$ws.Cells.Item(580, 6).Value = -159
$ws.Cells.Item(580, 7).Value = -2138.55
$ws.Cells.Item(583, 2).Value = 65066
$ws.Cells.Item(583, 5).Value = 13.61
$ws.Cells.Item(583, 6).Value = 313
$ws.Cells.Item(583, 7).Value = 4009.53
$ws.Cells.Item(584, 2).Value = 53263
$ws.Cells.Item(584, 5).Value = 15.29
$ws.Cells.Item(584, 6).Value = -309
$ws.Cells.Item(584, 7).Value = -3958.29
$ws.Cells.Item(599, 2).Value = 64925
$ws.Cells.Item(599, 5).Value = 13.97
$ws.Cells.Item(599, 6).Value = 302
$ws.Cells.Item(599, 7).Value = 3971.3
$ws.Cells.Item(600, 2).Value = 45709
$ws.Cells.Item(600, 5).Value = 15.69
$ws.Cells.Item(600, 6).Value = -300
$ws.Cells.Item(600, 7).Value = -3945
$ws.Cells.Item(604, 2).Value = 65067
$ws.Cells.Item(604, 5).Value = 15.65
$ws.Cells.Item(604, 6).Value = 338
$ws.Cells.Item(604, 7).Value = 4978.74
$ws.Cells.Item(605, 2).Value = 53595
$ws.Cells.Item(605, 5).Value = 17.61
$ws.Cells.Item(605, 6).Value = -335
$ws.Cells.Item(605, 7).Value = -4934.55
$ws.Cells.Item(720, 2).Value = 60022
$ws.Cells.Item(720, 5).Value = 37.22
$ws.Cells.Item(720, 6).Value = -113
$ws.Cells.Item(720, 7).Value = -3709.79
$ws.Cells.Item(721, 2).Value = 64830
$ws.Cells.Item(721, 5).Value = 34.9
$ws.Cells.Item(721, 6).Value = 117
$ws.Cells.Item(721, 7).Value = 3841.11
$ws.Cells.Item(872, 2).Value = 65079
$ws.Cells.Item(872, 6).Value = 21
$ws.Cells.Item(872, 7).Value = 858.27
$ws.Cells.Item(873, 2).Value = 65362
$ws.Cells.Item(873, 6).Value = 2
$ws.Cells.Item(873, 7).Value = 81.73999999999999
